{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find(\n  (p) => p.text.trim() === \"Jazykov\u00e9 prost\u0159edky\"\n);\nif (!target) {\n  throw new Error(\"Paragraph 'Jazykov\u00e9 prost\u0159edky' not found\");\n}\n\n// Append \" a Tropy\" right after the existing heading text, so the\n// heading reads \"Jazykov\u00e9 prost\u0159edky a Tropy\".\ntarget.getRange(\"End\").insertText(\" a Tropy\", \"End\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$found = $false\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    $text = $r.Text.TrimEnd([char]13, [char]7)\n    if ($text -eq \"Jazykov\u00e9 prost\u0159edky\") {\n        $r.InsertAfter(\" a Tropy\")\n        $found = $true\n        break\n    }\n}\nif (-not $found) {\n    throw \"Paragraph 'Jazykov\u00e9 prost\u0159edky' not found\"\n}\n"}
